# cambios finales de las plantillas y enlace
#
# Update the "STOCK MAXIMO" (V) / "STOCK MINIMO" (W) figures on the
# products sheet and leave the sheet scrolled/selected on the last
# cell the author ended up editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("products")

# Stock maximo / stock minimo values for each product row
$ws.Range("V2").Value = 5
$ws.Range("W2").Value = 15

$ws.Range("V3").Value = 10
$ws.Range("W3").Value = 22

$ws.Range("V4").Value = 15
$ws.Range("W4").Value = 33

$ws.Range("V5").Value = 15
$ws.Range("W5").Value = 33

# Leave the sheet active with the view scrolled further right and the
# final selection on W6, matching where the author left off editing.
$ws.Activate()
$win = $excel.Application.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 14
$ws.Range("W6").Select()
